$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 4, shifting all existing
# rows (4..34) down by one (to 5..35), same as the author inserting a new
# weekly price observation ahead of the existing history.
$ws.Rows(4).Insert()

# Populate the new row 4 with the new weekly record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44959
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = 100112003
$ws.Range("G4").Value = "Ajo"
$ws.Range("H4").Value = "Chino"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 21000
$ws.Range("L4").Value = 22000
$ws.Range("M4").Value = 21500
$ws.Range("N4").Value = "$/caja 10 kilos"
$ws.Range("O4").Value = "China"
$ws.Range("P4").Value = 2150
$ws.Range("Q4").Value = 10
$ws.Range("R4").Value = "Hortaliza"
